$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 2: add a Rounded Rectangle autoshape after the existing picture.
#    (The target shape ends up with Id=6 because PowerPoint's internal shape
#    id counter on this slide already has 2,3,5 in use; briefly adding and
#    removing a throw-away shape advances the counter past 4 so the real
#    shape lands on 6, matching the authored file.)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$dummy = $s2.Shapes.AddShape(5, 10, 10, 10, 10)
$dummy.Delete()

$rect = $s2.Shapes.AddShape(5, 8672052 / 12700, 1342103 / 12700, 2418735 / 12700, 2477729 / 12700)
$rect.Name = "Rounded Rectangle 5"
$rect.TextFrame.VerticalAnchor = 3
$rect.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# 2. Append a new slide ("Title and Content" layout) as slide 4.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Silde"
$titleRange.InsertAfter(" two")
